$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "galdiolo"
$ws.Range("A2").Value = "flor"
$ws.Range("A3").Value = "palmera"
$ws.Range("A4").Value = "bloso"
